$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 3080
$ws1.Range("F10").Value = 662
$ws1.Range("F14").Value = 2169
$ws1.Range("F16").Value = 754
$ws1.Range("F19").Value = 2692
$ws1.Range("F25").Value = 701
$ws1.Range("F26").Value = 701
$ws1.Range("F27").Value = 25
$ws1.Range("F35").Value = 916
$ws1.Range("F37").Value = 289

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 369
$ws2.Range("F15").Value = 279
$ws2.Range("F28").Value = 325
$ws2.Range("F38").Value = 617

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F6").Value = 302

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 302
$ws4.Range("F9").Value = 3080
$ws4.Range("F14").Value = 662
$ws4.Range("F18").Value = 369
$ws4.Range("F21").Value = 2169
$ws4.Range("F23").Value = 754
$ws4.Range("F27").Value = 2692
$ws4.Range("F29").Value = 279
$ws4.Range("F35").Value = 701
$ws4.Range("F36").Value = 701
$ws4.Range("F37").Value = 25
$ws4.Range("F46").Value = 916
$ws4.Range("F48").Value = 289
$ws4.Range("F50").Value = 617
$ws4.Range("F51").Value = 617
